# Apply symbol-list update (values scraped 2023-01-28 21:57 UTC).
# Cells are stored as text (inlineStr) in the workbook; numeric-looking
# values are written with a leading apostrophe so Excel keeps them as text
# instead of auto-converting to Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.79"
$ws.Range("E2").Value = "'-0.57%"
$ws.Range("D3").Value = "'38.80"
$ws.Range("E3").Value = "'6.36%"
$ws.Range("D4").Value = "'5.104"
$ws.Range("E4").Value = "'0.88%"
$ws.Range("D5").Value = "'0.08076"
$ws.Range("E5").Value = "'-0.33%"
$ws.Range("E6").Value = "'-4.90%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.200"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.009"
$ws.Range("E8").Value = "'1.94%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9270"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1437"
$ws.Range("E10").Value = "'-3.39%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1913"
$ws.Range("E11").Value = "'-1.50%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09151"
$ws.Range("E12").Value = "'0.74%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03505"
$ws.Range("E13").Value = "'-0.44%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09787"
$ws.Range("E14").Value = "'-1.01%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001395"
$ws.Range("E15").Value = "'-0.66%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005884"
$ws.Range("E16").Value = "'-3.93%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.779"
$ws.Range("E17").Value = "'-1.36%"
$ws.Range("D18").Value = "'3.344"
$ws.Range("E18").Value = "'-3.13%"
$ws.Range("E19").Value = "'0.27%"
$ws.Range("D20").Value = "'0.1294"
$ws.Range("E20").Value = "'-0.77%"
$ws.Range("D21").Value = "'4.672"
$ws.Range("E21").Value = "'-3.33%"
$ws.Range("E22").Value = "'3.12%"
$ws.Range("D23").Value = "'0.04380"
$ws.Range("E23").Value = "'-0.13%"
$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-0.57%"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("E25").Value = "'1.97%"
$ws.Range("E26").Value = "'-0.05%"
$ws.Range("D39").Value = "'0.02036"
$ws.Range("E39").Value = "'-0.87%"
$ws.Range("E40").Value = "'-1.51%"
$ws.Range("D41").Value = "'0.007525"
$ws.Range("E41").Value = "'0.38%"
$ws.Range("D42").Value = "'0.009712"
$ws.Range("E42").Value = "'-3.38%"
$ws.Range("D43").Value = "'0.1341"
$ws.Range("E43").Value = "'-2.08%"
$ws.Range("D44").Value = "'0.002104"
$ws.Range("E44").Value = "'-0.99%"
$ws.Range("D45").Value = "'0.009901"
$ws.Range("E45").Value = "'-0.15%"
$ws.Range("D46").Value = "'0.00006215"
$ws.Range("E47").Value = "'-0.02%"
$ws.Range("E49").Value = "'12.57%"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'-0.02%"
